$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 58, shifting existing rows 58:90 down to 59:91
$ws.Rows(58).Insert()

# Populate the newly inserted row 58 with the new weekly data point.
# Static/metadata columns mirror the rest of this market/category block.
$ws.Range("A58").Value = 6
$ws.Range("B58").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C58").Value = "Metropolitana"
$ws.Range("D58").Value = 44762
$ws.Range("E58").Value = 13
$ws.Range("F58").Value = 100114007
$ws.Range("G58").Value = "Jengibre"
$ws.Range("H58").Value = "Sin especificar"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 310
$ws.Range("K58").Value = 12000
$ws.Range("L58").Value = 13000
$ws.Range("M58").Value = 12581
$ws.Range("N58").Value = "$/caja 13 kilos"
$ws.Range("O58").Value = "Perú"
$ws.Range("P58").Value = 968
$ws.Range("Q58").Value = 13
$ws.Range("R58").Value = "Hortaliza"
